# Add "Ensemble" model row (row 6) to the training metrics table, mirroring
# the existing data rows (2-5): a bold/centered/bordered label cell in
# column A followed by six numeric metric values in columns B-G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Label cell, formatted like the other model-name cells (A2:A5).
$ws.Range("A6").Value = "Ensemble"
$labelCell = $ws.Range("A6")
$labelCell.Font.Bold = $true
$labelCell.HorizontalAlignment = -4108   # xlCenter
$labelCell.VerticalAlignment = -4160     # xlTop
$labelCell.Borders.LineStyle = 1         # xlContinuous (thin box border)

# Metric values for the new "Ensemble" row.
$ws.Range("B6").Value = 0.578
$ws.Range("C6").Value = 0.543
$ws.Range("D6").Value = 0.46
$ws.Range("E6").Value = 0.678
$ws.Range("F6").Value = 0.499
$ws.Range("G6").Value = 0.875
